{"js": "// Problem 2 solution evaluation\n// After the \"Identify potential solutions\" answer paragraph (the one ending\n// \"...required amount to absolutely guarantee a solution every time.\"),\n// insert three new paragraphs:\n//   1) a blank spacer paragraph\n//   2) \"Evaluate each potential solution:\"\n//   3) the evaluation text for Problem 2.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker =\n  \"required amount to absolutely guarantee a solution every time\";\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf(marker) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the Problem 2 'potential solutions' paragraph.\");\n}\n\n// Insert the blank spacer paragraph right after the target paragraph.\nconst blank = target.insertParagraph(\"\", \"After\");\n\n// Insert the \"Evaluate each potential solution:\" heading paragraph.\nconst heading = blank.insertParagraph(\"Evaluate each potential solution:\", \"After\");\n\n// Insert the evaluation paragraph.\nheading.insertParagraph(\n  \"Each of these solutions will meet the goals.  Even though you have a possibility of solving the problems with less socks, these solutions guarantee the problem is solved because it looks at the worst case scenario.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Problem 2 solution evaluation\n# After the \"Identify potential solutions\" answer paragraph (the one ending\n# \"...required amount to absolutely guarantee a solution every time.\"),\n# insert three new paragraphs:\n#   1) a blank spacer paragraph\n#   2) \"Evaluate each potential solution:\"\n#   3) the evaluation text for Problem 2.\n\n$d = $word.ActiveDocument\n\n$marker = \"required amount to absolutely guarantee a solution every time\"\n\n$targetIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*$marker*\") {\n        $targetIdx = $i\n        break\n    }\n}\n\nif ($targetIdx -eq -1) {\n    throw \"Could not find the Problem 2 'potential solutions' paragraph.\"\n}\n\n# Insert a blank spacer paragraph right after the target paragraph.\n$target = $d.Paragraphs.Item($targetIdx)\n$target.Range.InsertParagraphAfter()\n\n# Insert the \"Evaluate each potential solution:\" heading paragraph.\n$blank = $d.Paragraphs.Item($targetIdx + 1)\n$blank.Range.InsertParagraphAfter()\n\n$heading = $d.Paragraphs.Item($targetIdx + 2)\n$heading.Range.Text = \"Evaluate each potential solution:\"\n\n# Insert the evaluation paragraph.\n$heading.Range.InsertParagraphAfter()\n\n$evalPara = $d.Paragraphs.Item($targetIdx + 3)\n$evalPara.Range.Text = \"Each of these solutions will meet the goals.  Even though you have a possibility of solving the problems with less socks, these solutions guarantee the problem is solved because it looks at the worst case scenario.\"\n"}
